$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.07147868145859068
$ws.Cells.Item(2, 8).Value = 114.6830761336013
$ws.Cells.Item(2, 9).Value = -47.92863797341975
$ws.Cells.Item(3, 7).Value = 0.0449230761555611
$ws.Cells.Item(3, 8).Value = -10.98118480765538
$ws.Cells.Item(4, 7).Value = 0.03026286798604583
$ws.Cells.Item(4, 8).Value = -37.82772071640407
$ws.Cells.Item(5, 7).Value = 0.02540116381037075
$ws.Cells.Item(5, 8).Value = -64.04879019512209
$ws.Cells.Item(6, 7).Value = -0.07779056171429792
$ws.Cells.Item(6, 8).Value = 33.86600407384236
$ws.Cells.Item(7, 7).Value = -0.07254255638893733
$ws.Cells.Item(7, 8).Value = 41.98149007573591
$ws.Cells.Item(8, 7).Value = -0.295217284720253
$ws.Cells.Item(8, 8).Value = -48.21727915611579
$ws.Cells.Item(9, 7).Value = -0.221924330262072
$ws.Cells.Item(9, 8).Value = 26.53382642777011
$ws.Cells.Item(10, 7).Value = -0.03342811851076076
$ws.Cells.Item(10, 8).Value = -2726.520277607825
$ws.Cells.Item(11, 7).Value = 0.04347243055808837
$ws.Cells.Item(11, 8).Value = 311.4069819073134
$ws.Cells.Item(12, 7).Value = 0.1837474491125869
$ws.Cells.Item(12, 8).Value = -13.29038684079105
$ws.Cells.Item(13, 7).Value = 0.23250959624587
$ws.Cells.Item(13, 8).Value = -1.058310315862158
$ws.Cells.Item(14, 7).Value = -0.06580521124171784
$ws.Cells.Item(14, 8).Value = 27.72650567187113
$ws.Cells.Item(15, 7).Value = -0.06715710993832892
$ws.Cells.Item(15, 8).Value = 5.370847259901038
$ws.Cells.Item(16, 7).Value = 0.1907936063235357
$ws.Cells.Item(16, 8).Value = -0.3445606109700473
$ws.Cells.Item(17, 7).Value = 0.1696943861460263
$ws.Cells.Item(17, 8).Value = -2.401360900785445
$ws.Cells.Item(18, 7).Value = 0.02260558681381364
$ws.Cells.Item(18, 8).Value = -58.31523744577826
$ws.Cells.Item(19, 7).Value = 0.07478950884487799
$ws.Cells.Item(19, 8).Value = -13.07548184839422
$ws.Cells.Item(20, 7).Value = 0.007373700479570815
$ws.Cells.Item(20, 8).Value = -42.0712236218401
$ws.Cells.Item(21, 7).Value = -0.03725811028716507
$ws.Cells.Item(21, 8).Value = 30.84186070793264
$ws.Cells.Item(22, 7).Value = 0.06924214092562987
$ws.Cells.Item(22, 8).Value = 6.081253621223612
$ws.Cells.Item(23, 7).Value = 0.07765514047267943
$ws.Cells.Item(23, 8).Value = 34.64639722719058
$ws.Cells.Item(24, 7).Value = 0.03790848929130086
$ws.Cells.Item(24, 8).Value = 17.02278274351985
$ws.Cells.Item(25, 7).Value = 0.04305493236069727
$ws.Cells.Item(25, 8).Value = 46.28547708903096
$ws.Cells.Item(26, 7).Value = 0.1073982701158028
$ws.Cells.Item(26, 8).Value = -5.201757105776029
$ws.Cells.Item(27, 7).Value = 0.1398827576158806
$ws.Cells.Item(27, 8).Value = 55.10188594234574
$ws.Cells.Item(28, 7).Value = 0.1131630080145592
$ws.Cells.Item(28, 8).Value = -3.664096775921789
$ws.Cells.Item(29, 7).Value = 0.1439690968699943
$ws.Cells.Item(29, 8).Value = 20.34766160512822
$ws.Cells.Item(30, 7).Value = 0.07030092853317786
$ws.Cells.Item(30, 8).Value = 4.570386923710991
$ws.Cells.Item(31, 7).Value = 0.07086164524611151
$ws.Cells.Item(31, 8).Value = 3.270270529709203
$ws.Cells.Item(32, 7).Value = 0.05320511211118779
$ws.Cells.Item(32, 8).Value = 21.84407710166562
$ws.Cells.Item(33, 7).Value = 0.0482339694822452
$ws.Cells.Item(33, 8).Value = -11.23411723547709
$ws.Cells.Item(34, 7).Value = -0.007082817274089428
$ws.Cells.Item(34, 8).Value = 62.91450447369809
$ws.Cells.Item(35, 7).Value = 0.01344832753129484
$ws.Cells.Item(35, 8).Value = -3.614656444936093
$ws.Cells.Item(36, 7).Value = -0.00788876241927247
$ws.Cells.Item(36, 8).Value = -151.0304970982246
$ws.Cells.Item(37, 7).Value = -0.006905272940685193
$ws.Cells.Item(37, 8).Value = -155.1377681775152
$ws.Cells.Item(38, 7).Value = 0.0372065149879895
$ws.Cells.Item(38, 8).Value = -48.13902560123869
$ws.Cells.Item(39, 7).Value = 0.0663204996857413
$ws.Cells.Item(39, 8).Value = 54.00984737733414
$ws.Cells.Item(40, 7).Value = 0.03232882613625514
$ws.Cells.Item(40, 8).Value = -27.72413392096107
$ws.Cells.Item(41, 7).Value = 0.0296526909029519
$ws.Cells.Item(41, 8).Value = 139.9216585444044
$ws.Cells.Item(42, 7).Value = 0.06117008341610763
$ws.Cells.Item(42, 8).Value = 17.0070367501139
$ws.Cells.Item(43, 7).Value = 0.07871099396561926
$ws.Cells.Item(43, 8).Value = 57.73782546958969
$ws.Cells.Item(44, 7).Value = 0.1215498708675968
$ws.Cells.Item(44, 8).Value = -7.742639907706141
$ws.Cells.Item(45, 7).Value = 0.1546971482643252
$ws.Cells.Item(45, 8).Value = -13.78253063288607
$ws.Cells.Item(46, 7).Value = -0.05368190259391904
$ws.Cells.Item(46, 8).Value = -22.19917057921221
$ws.Cells.Item(47, 7).Value = -0.01169006436056031
$ws.Cells.Item(47, 8).Value = -346.2594503754428
$ws.Cells.Item(48, 7).Value = 0.01187591370294495
$ws.Cells.Item(48, 8).Value = -18.05092055492683
$ws.Cells.Item(49, 7).Value = -0.001308095717333042
$ws.Cells.Item(49, 8).Value = 76.46809781526318
$ws.Cells.Item(50, 7).Value = 0.09705183841290542
$ws.Cells.Item(50, 8).Value = -32.09904159927049
$ws.Cells.Item(51, 7).Value = 0.1465823103934866
$ws.Cells.Item(51, 8).Value = 11.92461360946671
$ws.Cells.Item(52, 7).Value = 0.0561697084943592
$ws.Cells.Item(52, 8).Value = -9.334533309351031
$ws.Cells.Item(53, 7).Value = 0.04773645594918585
$ws.Cells.Item(53, 8).Value = -21.97583294506032
$ws.Cells.Item(54, 7).Value = -0.1433159341703481
$ws.Cells.Item(54, 8).Value = -60.47569697225157
$ws.Cells.Item(55, 7).Value = -0.09126890292579742
$ws.Cells.Item(55, 8).Value = 12.06100425614045
$ws.Cells.Item(56, 7).Value = 0.1632268906029936
$ws.Cells.Item(56, 8).Value = 5.223747139631691
$ws.Cells.Item(57, 7).Value = 0.1854674612244922
$ws.Cells.Item(57, 8).Value = 32.96483631506072
